# Content from David, tweaks from Domenic

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Content tweak: fix the Bottou SGD Tricks link text (D4) ---
# Original: <a href="refs/bottou-sgd-tricks-2012.pdf>Bottou's SGD Tricks</a>
# New:      <a href="refs/bottou-sgd-tricks-2012.pdf>Bottou's SGD Tricks"</a>
$ws.Range("D4").Value = '<a href="refs/bottou-sgd-tricks-2012.pdf>Bottou''s SGD Tricks"</a>'

# D5 text stays "BV Preface, Ch 1" (unchanged content, just reordered in the
# shared-string table as a side effect of the edit above).
$ws.Range("D5").Value = "BV Preface, Ch 1"

# --- View tweaks from Domenic ---
# Move the active selection to A11 (also resets the scrolled-to-B1 view,
# restoring column A into view).
$ws.Range("A11").Select() | Out-Null
